$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 12, shifting existing rows 12..79 down to 13..80.
$ws.Rows.Item(12).Insert()

# Populate the newly inserted (blank) row 12 with the new data record.
$ws.Cells.Item(12, 1).Value = 3
$ws.Cells.Item(12, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(12, 3).Value = "Coquimbo"
$ws.Cells.Item(12, 4).Value = 44802
$ws.Cells.Item(12, 5).Value = 5
$ws.Cells.Item(12, 6).Value = 100112035
$ws.Cells.Item(12, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 73
$ws.Cells.Item(12, 11).Value = 15000
$ws.Cells.Item(12, 12).Value = 16000
$ws.Cells.Item(12, 13).Value = 15479
$ws.Cells.Item(12, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(12, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(12, 16).Value = 1032
$ws.Cells.Item(12, 17).Value = 15
$ws.Cells.Item(12, 18).Value = "Hortaliza"
